# Fruta / hortaliza, semanal
# Re-sort the data rows (2..14) by the "Fecha" (date, column D) ascending,
# moving D, J, K, L, M, O, P along with each record.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values per row (2..14) for columns: D (Fecha), J (Volumen),
# K (Precio minimo), L (Precio maximo), M (Precio promedio ponderado),
# O (Origen), P (Precio $/Kg)
$rows = @(
    @{ Row = 2;  D = 44432; J = 25; K = 14000; L = 14000; M = 14000; O = "Provincia del Elquí"; P = 467 },
    @{ Row = 3;  D = 44460; J = 45; K = 13000; L = 13000; M = 13000; O = "Provincia de Limarí";  P = 433 },
    @{ Row = 4;  D = 44467; J = 35; K = 12000; L = 12000; M = 12000; O = "Provincia de Limarí";  P = 400 },
    @{ Row = 5;  D = 44435; J = 25; K = 14000; L = 14000; M = 14000; O = "Provincia de Limarí";  P = 467 },
    @{ Row = 6;  D = 44435; J = 25; K = 14000; L = 14000; M = 14000; O = "Provincia del Elquí"; P = 467 },
    @{ Row = 7;  D = 44446; J = 25; K = 14000; L = 14000; M = 14000; O = "Provincia de Limarí";  P = 467 },
    @{ Row = 8;  D = 44474; J = 45; K = 10000; L = 10000; M = 10000; O = "Provincia de Limarí";  P = 333 },
    @{ Row = 9;  D = 44418; J = 30; K = 15000; L = 15000; M = 15000; O = "Provincia de Limarí";  P = 500 },
    @{ Row = 10; D = 44376; J = 25; K = 18000; L = 18000; M = 18000; O = "Provincia de Limarí";  P = 600 },
    @{ Row = 11; D = 44453; J = 50; K = 12000; L = 12000; M = 12000; O = "Provincia de Limarí";  P = 400 },
    @{ Row = 12; D = 44425; J = 35; K = 14000; L = 14000; M = 14000; O = "Provincia de Limarí";  P = 467 },
    @{ Row = 13; D = 44449; J = 45; K = 12000; L = 12000; M = 12000; O = "Provincia de Limarí";  P = 400 },
    @{ Row = 14; D = 44421; J = 25; K = 15000; L = 16000; M = 15400; O = "Provincia de Limarí";  P = 513 }
)

foreach ($r in $rows) {
    $ws.Cells.Item($r.Row, 4).Value  = $r.D   # D: Fecha
    $ws.Cells.Item($r.Row, 10).Value = $r.J   # J: Volumen
    $ws.Cells.Item($r.Row, 11).Value = $r.K   # K: Precio minimo
    $ws.Cells.Item($r.Row, 12).Value = $r.L   # L: Precio maximo
    $ws.Cells.Item($r.Row, 13).Value = $r.M   # M: Precio promedio ponderado
    $ws.Cells.Item($r.Row, 15).Value = $r.O   # O: Origen
    $ws.Cells.Item($r.Row, 16).Value = $r.P   # P: Precio $/Kg
}
